$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.040.49"
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = "'2.412.98"
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'553.79"
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = "'136.49"
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.77%  '
$ws.Range('E9').Value = '  -0.86%  '
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('D13').Value = "'24.75"
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').Value = "'2.844.73"
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').Value = "'59.958.36"
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').Value = "'2.392.98"
$ws.Range('E17').Value = '  -1.38%  '
$ws.Range('D18').Value = "'11.20"
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('E19').Value = '  +3.45%  '
$ws.Range('D20').Value = "'326.61"
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = "'64.74"
$ws.Range('E23').Value = '  -1.52%  '
$ws.Range('E24').Value = '  +5.57%  '
$ws.Range('D25').Value = "'8.62"
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E27').Value = '  +4.48%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = "'1.78"
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = "'0.0₃0772"
$ws.Range('E29').Value = '  -1.29%  '
$ws.Range('D30').Value = "'170.87"
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('E32').Value = '  +5.68%  '
$ws.Range('E33').Value = '  -3.34%  '
$ws.Range('D34').Value = "'18.41"
$ws.Range('E34').Value = '  -1.03%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  +2.16%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('D39').Value = "'324.08"
$ws.Range('E39').Value = '  +3.21%  '
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = "'146.49"
$ws.Range('E41').Value = '  +5.09%  '
$ws.Range('D42').Value = "'3.61"
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('D44').Value = "'19.77"
$ws.Range('E44').Value = '  +2.41%  '
$ws.Range('D45').Value = "'0.0515"
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').Value = "'0.578"
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('D51').Value = "'0.939"
$ws.Range('E51').Value = '  -1.73%  '
